$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) after the existing "sum" column (G).
# Copy the header style from the existing header cell (G1) so the new
# header (H1) gets the same bold/centered/bordered formatting, then
# overwrite its text with "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Fill in the Save column values for the data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
